# [Kadastro App] Yeni kayit eklendi: 2985
# Appends the new record (row 51) to both the master "Kayitlar" sheet and
# the filtered "Erdemli" sheet, which mirror each other for this record's
# Birim ("Erdemli").

$wb = $excel.ActiveWorkbook

$newRow = @("'2985", "'2025-09-10", "'Erdemli", "'1", "'LİHKAB", "'SEVİL SARAÇER (Tekniker), ÖZKAN AKBAŞ (Mühendis)")

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item(51, $col)
        $cell.Value = $newRow[$col - 1]
        $cell.Style = "Normal"
    }
}
